$d = $word.ActiveDocument

# --- 1. Remove the stray _GoBack bookmark that sat on the
#        "IMPORTANT BELOW, FIX NOW" paragraph. ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- 2. Append a blank paragraph followed by a new bold paragraph
#        (with the _GoBack bookmark now wrapping its run) right
#        before the final section break. ---
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

$blankParaXml = "<w:p $wNs/>"

$newParaXml = "<w:p $wNs>" +
    "<w:pPr><w:rPr><w:b/></w:rPr></w:pPr>" +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    "<w:r><w:rPr><w:b/></w:rPr>" +
    "<w:t>These issues were just from one major bug, which has since been resolved.</w:t>" +
    "</w:r>" +
    '<w:bookmarkEnd w:id="0"/>' +
    "</w:p>"

$insertPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$null = $insertPoint.InsertXML($blankParaXml + $newParaXml)
